# Fruta / hortaliza, semanal
# Insert a new data row at row 111 (shifting the existing rows 111-176 down
# to 112-177) and populate the new row with the latest weekly price entry
# for "Ajo" (Terminal Hortofruticola Agro Chillan).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push every existing row from 111 down by one to make room for the new entry.
$ws.Rows(111).Insert()

$newRow = 111

$ws.Cells.Item($newRow, 1).Value2  = 7
$ws.Cells.Item($newRow, 2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($newRow, 3).Value2  = "Ñuble"
$ws.Cells.Item($newRow, 4).Value2  = 44567
$ws.Cells.Item($newRow, 5).Value2  = 16
$ws.Cells.Item($newRow, 6).Value2  = 100112003
$ws.Cells.Item($newRow, 7).Value2  = "Ajo"
$ws.Cells.Item($newRow, 8).Value2  = "Chino"
$ws.Cells.Item($newRow, 9).Value2  = "Primera"
$ws.Cells.Item($newRow, 10).Value2 = 80
$ws.Cells.Item($newRow, 11).Value2 = 19000
$ws.Cells.Item($newRow, 12).Value2 = 20000
$ws.Cells.Item($newRow, 13).Value2 = 19500
$ws.Cells.Item($newRow, 14).Value2 = "`$/caja 10 kilos"
$ws.Cells.Item($newRow, 15).Value2 = "China"
$ws.Cells.Item($newRow, 16).Value2 = 1950
$ws.Cells.Item($newRow, 17).Value2 = 10
$ws.Cells.Item($newRow, 18).Value2 = "Hortaliza"
